$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1
$ws.Range("A1").Value = 45436

# TIRAFONDOS DE 3/16 (5MM)
$ws.Range("D32").Value = 2236.3
$ws.Range("D33").Value = 2673
$ws.Range("D34").Value = 2992
$ws.Range("D35").Value = 3223
$ws.Range("D36").Value = 3410
$ws.Range("D37").Value = 3829.1
$ws.Range("D38").Value = 4213
$ws.Range("D39").Value = 4510
$ws.Range("D40").Value = 4752

# TIRAFONDOS DE 1/4 (6,5MM)
$ws.Range("D46").Value = 2728
$ws.Range("D47").Value = 3047
$ws.Range("D48").Value = 3344
$ws.Range("D49").Value = 3894
$ws.Range("D50").Value = 4236.1
$ws.Range("D51").Value = 4719
$ws.Range("D52").Value = 5297.6
$ws.Range("D53").Value = 5522
$ws.Range("D54").Value = 6094
$ws.Range("D55").Value = 6809
$ws.Range("D56").Value = 7393.1
$ws.Range("D57").Value = 8305
$ws.Range("D58").Value = 9086
$ws.Range("D59").Value = 9955
$ws.Range("D60").Value = 11348.823
$ws.Range("D61").Value = 12044.74

# TIRAFONDOS DE 5/16 (8MM)
$ws.Range("D67").Value = 5049
$ws.Range("D68").Value = 5758.5
$ws.Range("D69").Value = 6204
$ws.Range("D70").Value = 6303
$ws.Range("D71").Value = 7117
$ws.Range("D72").Value = 7667
$ws.Range("D73").Value = 8538.2
$ws.Range("D74").Value = 9168.5
$ws.Range("D75").Value = 10642.5
$ws.Range("D76").Value = 11572
$ws.Range("D77").Value = 12881
$ws.Range("D78").Value = 14179
$ws.Range("D79").Value = 15730
$ws.Range("D80").Value = 21703
$ws.Range("D81").Value = 24178

# TIRAFONDOS DE 3/8 (9,5MM)
$ws.Range("D87").Value = 8954
$ws.Range("D88").Value = 9240
$ws.Range("D89").Value = 10824
$ws.Range("D90").Value = 11858
$ws.Range("D91").Value = 13090
$ws.Range("D92").Value = 15994
$ws.Range("D93").Value = 16643
$ws.Range("D94").Value = 18480
$ws.Range("D95").Value = 20240
$ws.Range("D96").Value = 22385
